# Apply AutoFilter on Table1 so that only rows where BlackPlayer = "Radjabov"
# remain visible (this hides all other data rows), then move the active
# cell/selection on Sheet1 to B73.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lo = $ws.ListObjects.Item(1)

# Column 3 of the table range (A:F) is "BlackPlayer" (table column id 3,
# i.e. 0-indexed colId 2 in the OOXML <filterColumn>), filtered to the
# single value "Radjabov" using the standard values filter (xlFilterValues = 7).
$lo.Range.AutoFilter(3, @("Radjabov"), 7)

# Update the selection/active cell shown when the workbook is reopened.
$ws.Range("B73").Select()
